$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 39 (old row 39 "Structures de boucles" and
# the rows after it shift down by one).
$ws.Rows.Item(39).Insert()

# Old row 42 ("fin") is now at row 43. Insert 4 new rows before it, for the
# new "plot" / "Boucles IF FOR WHILE" / "algo" exercises.
$ws.Range("A43:A46").EntireRow.Insert()

# --- Fill the new content (order matches the original authoring session) ---

# New "plot" exercises (rows 44-45).
$ws.Cells.Item(44, 1).Value = "plot"
$ws.Cells.Item(44, 2).Value = "PLT-002"
$ws.Cells.Item(44, 3).Value = "Tracé de courbes"

$ws.Cells.Item(45, 1).Value = "plot"
$ws.Cells.Item(45, 2).Value = "PLT-003"
$ws.Cells.Item(45, 3).Value = "Exemple avancé"

# New row 39: "Fonctions en Python" exercise.
$ws.Cells.Item(39, 2).Value = "PYB-308"
$ws.Cells.Item(39, 3).Value = "Fonctions en Python"
$ws.Cells.Item(39, 1).Value = "python_bases"
$ws.Cells.Item(39, 5).Value = "301+305"

# New row 43: "Boucles IF, FOR, WHILE" exercise.
# Leading apostrophe forces this comma-grouped digit string to stay text
# instead of being auto-parsed as the number 400401404; reapplying the
# "Normal" style afterwards drops the quote-prefix formatting flag again.
$ws.Cells.Item(43, 5).Value = "'400,401,404"
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(43, 1).Value = "python_bases"
$ws.Cells.Item(43, 2).Value = "PYB-411"
$ws.Cells.Item(43, 3).Value = "Boucles IF, FOR, WHILE"

# New row 46: "algo" exercise.
$ws.Cells.Item(46, 3).Value = "Algorithme glouton -- Problème du rendu de monnaie"
$ws.Cells.Item(46, 2).Value = "ALGO-012"
$ws.Cells.Item(46, 1).Value = "algo"

# Update the view so the last edited rows are visible, matching the
# author's final scroll/selection position.
$ws.Range("A47").Select()
$excel.ActiveWindow.ScrollRow = 19
